# ---------------------------------------------------------------------------
# Edit: Slide 1 subtitle credits line gains a third author ("Dmitry
# Fedorovichev") after "Kirill Gelvan, Mikhail Kuznetsov", and two table
# cells on Slide 10 get a "<" prefix.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 1 - add "Dmitry Fedorovichev" to the credits paragraph.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# "Subtitle 2" shape holding the two-line author credit block is shape 5.
$creditShape = $s1.Shapes.Item(5)
$tr = $creditShape.TextFrame.TextRange

# Second paragraph currently reads "Kirill Gelvan, Mikhail Kuznetsov" followed
# by a superscript "1" footnote marker run.
$para2 = $tr.Paragraphs(2)

# Split "Kuznetsov" (characters 24-32 of the paragraph) into its own run so
# the new content can be inserted right after it, ahead of the superscript.
$kuznetsovRange = $para2.Characters(24, 9)
$kuznetsovRange.Text = $kuznetsovRange.Text

# Insert the new names right after "Kuznetsov" (and before the superscript
# "1" run, which stays untouched at the end of the paragraph).
$kuznetsovRange2 = $para2.Characters(24, 9)
[void]$kuznetsovRange2.InsertAfter(", Dmitry Fedorovichev")

# Break the newly-inserted text into its own runs: ", ", "Dmitry " and
# "Fedorovichev".
$commaRange = $para2.Characters(33, 2)
$commaRange.Text = $commaRange.Text

$dmitryRange = $para2.Characters(35, 7)
$dmitryRange.Text = $dmitryRange.Text
$dmitryRange.Font.Name = ""

$fedorovichevRange = $para2.Characters(42, 12)
$fedorovichevRange.Text = $fedorovichevRange.Text
$fedorovichevRange.Font.Name = ""

# ---------------------------------------------------------------------------
# 2) Slide 10 - prefix two summary-row table cells with "<".
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$tableShape = $s10.Shapes.Item(4)
$tbl = $tableShape.Table

$sumInstancesCell = $tbl.Cell(11, 2)
$sumInstancesCell.Shape.TextFrame.TextRange.Text = "<194``974"

$sumPercentCell = $tbl.Cell(11, 3)
$sumPercentCell.Shape.TextFrame.TextRange.Text = "<86.44"
